$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued columns (Coin name / Link) - plain string assignment is fine since
# these values are not numeric-looking, so Excel keeps them as text automatically.
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

# Numeric-looking text columns (Price / Volume%) - force Text number format first so
# Excel stores the literal digit string instead of auto-converting to a Double/percentage,
# which would lose formatting (e.g. "-0.66%" must stay literal text, not become -0.0066).
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value2 = '245.16'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value2 = '-0.66%'
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value2 = '28.68'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value2 = '-2.96%'
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value2 = '5.254'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value2 = '1.59%'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = '0.05708'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value2 = '0.05%'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = '6.614'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value2 = '0.18%'
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = '3.178'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value2 = '3.19%'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value2 = '0.8532'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value2 = '-0.50%'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = '0.8600'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value2 = '-1.86%'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = '0.1373'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value2 = '0.29%'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value2 = '0.07044'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value2 = '-0.46%'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = '0.03161'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value2 = '10.14%'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = '0.09292'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value2 = '-0.97%'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = '0.001528'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value2 = '1.03%'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = '0.0005950'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value2 = '-1.11%'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value2 = '0.005920'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value2 = '-1.85%'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = '3.490'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value2 = '0.14%'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value2 = '2.174'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value2 = '-4.57%'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = '0.3165'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value2 = '-0.30%'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = '0.03330'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value2 = '0.64%'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value2 = '-1.78%'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value2 = '3.493'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value2 = '0.80%'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = '0.04112'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value2 = '-1.57%'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = '0.1379'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value2 = '-0.04%'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = '0.001219'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value2 = '-0.04%'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value2 = '0.004143'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value2 = '-17.83%'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value2 = '-0.77%'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = '0.0001449'
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value2 = '-25.25%'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = '0.03766'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value2 = '0.36%'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = '0.1065'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value2 = '-0.58%'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = '0.003689'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value2 = '-35.82%'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = '0.002450'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value2 = '16.69%'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = '0.009336'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value2 = '-8.61%'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value2 = '0.00005295'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value2 = '2.67%'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = '0.00000000750'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value2 = '0.05%'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = '0.07501'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value2 = '7.20%'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = '0.002431'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value2 = '-5.54%'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = '0.00002100'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value2 = '0.05%'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = '0.0002000'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value2 = '0.05%'
